# Adding Printer script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Printer" form-field label cell next to the existing labels
$cell = $ws.Range("F4")
$cell.Value2 = "Printer"

# Style it like the other field-type labels: bold, small "code" font, accent blue
$fnt = $cell.Font
$fnt.Name = "JetBrains Mono"
$fnt.Family = 3
$fnt.Bold = $true
$fnt.Size = 9.8
$fnt.Color = 14580521

$cell.VerticalAlignment = -4108

# Move the selection to the newly added cell
$cell.Select() | Out-Null

# Set the sheet to print in portrait orientation
$ws.PageSetup.Orientation = 1
